# Recommendation_Likelihood.xlsx — "Add files via upload"
#
# The author filled in the previously-blank "Survey 3" row (row 4) on
# Sheet1 with its response counts across the five Likelihood columns
# (B:F — Very Likely, Likely, May Recommend, Unlikely, Very Unlikely),
# and left the selection on that newly entered range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("Survey 3") was missing B:F values — fill them in.
$ws.Range("B4").Value = 22
$ws.Range("C4").Value = 22
$ws.Range("D4").Value = 31
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1

# Leave the active selection on the range that was just edited.
$ws.Range("B4:F4").Select()
